# Updated cryptos list on Sun Nov 19 23:19:58 UTC 2023 with GitHub Actions
# Refreshes the Price (D) / Volume(1h) (E) columns for every coin row, and
# for rows 49-50 also updates Coin (B) / Link (C) because MultiversX and
# MXToken swapped rank positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All of the source data lives in cells typed as inline text (t="inlineStr"),
# including values that look numeric (e.g. "1.00", "246.51"). Plain
# `Range.Value = ...` lets Excel auto-coerce such strings into numbers,
# which both changes the stored type and can silently drop formatting such
# as trailing zeros (e.g. "1.00" -> 1). Forcing the cell to Text format
# before the write keeps it a string; clearing formats afterwards drops the
# temporary "@" number format again so no stray style is left behind.
function Set-CellText($addr, [string]$val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

# Row 2
Set-CellText "D2" "37.430.93"
Set-CellText "E2" "  +2.34%  "

# Row 3
Set-CellText "D3" "2.000.38"
Set-CellText "E3" "  +2.09%  "

# Row 4
Set-CellText "E4" "  -0.02%  "

# Row 5
Set-CellText "D5" "246.51"
Set-CellText "E5" "  +0.80%  "

# Row 6
Set-CellText "D6" "0.632"
Set-CellText "E6" "  +2.40%  "

# Row 7
Set-CellText "D7" "62.06"
Set-CellText "E7" "  +5.65%  "

# Row 8
Set-CellText "D8" "1.00"
Set-CellText "E8" "  -0.03%  "

# Row 9
Set-CellText "D9" "0.386"
Set-CellText "E9" "  +2.06%  "

# Row 10
Set-CellText "D10" "0.0807"
Set-CellText "E10" "  -0.14%  "

# Row 11
Set-CellText "E11" "  +0.43%  "

# Row 12
Set-CellText "E12" "  +8.74%  "

# Row 13
Set-CellText "D13" "22.67"
Set-CellText "E13" "  +2.55%  "

# Row 14
Set-CellText "D14" "0.851"
Set-CellText "E14" "  +2.42%  "

# Row 15
Set-CellText "D15" "2.287.71"
Set-CellText "E15" "  +1.84%  "

# Row 16
Set-CellText "D16" "5.46"
Set-CellText "E16" "  +3.19%  "

# Row 17
Set-CellText "D17" "2.001.66"
Set-CellText "E17" "  +1.98%  "

# Row 18
Set-CellText "D18" "37.300.93"
Set-CellText "E18" "  +2.15%  "

# Row 19
Set-CellText "D19" "70.39"
Set-CellText "E19" "  +0.94%  "

# Row 20
Set-CellText "D20" "0.0₃0866"
Set-CellText "E20" "  +1.32%  "

# Row 21
Set-CellText "D21" "5.21"
Set-CellText "E21" "  +3.21%  "

# Row 22
Set-CellText "D22" "231.20"
Set-CellText "E22" "  +1.15%  "

# Row 23
Set-CellText "E23" "  +0.20%  "

# Row 24
Set-CellText "D24" "2.54"
Set-CellText "E24" "  +3.34%  "

# Row 25
Set-CellText "D25" "2.38"
Set-CellText "E25" "  +0.99%  "

# Row 26
Set-CellText "D26" "0.147"
Set-CellText "E26" "  +6.51%  "

# Row 27
Set-CellText "D27" "9.36"
Set-CellText "E27" "  +1.09%  "

# Row 28
Set-CellText "D28" "163.95"
Set-CellText "E28" "  +2.19%  "

# Row 29
Set-CellText "D29" "19.77"
Set-CellText "E29" "  +1.61%  "

# Row 30
Set-CellText "E30" "  +18.64%  "

# Row 31
Set-CellText "E31" "  +1.68%  "

# Row 32
Set-CellText "D32" "4.88"
Set-CellText "E32" "  +3.69%  "

# Row 33
Set-CellText "D33" "0.0626"
Set-CellText "E33" "  +1.05%  "

# Row 34
Set-CellText "D34" "4.62"
Set-CellText "E34" "  +6.83%  "

# Row 35
Set-CellText "D35" "2.33"
Set-CellText "E35" "  +4.05%  "

# Row 36
Set-CellText "D36" "1.00"
Set-CellText "E36" "  -0.10%  "

# Row 37
Set-CellText "D37" "3.37"
Set-CellText "E37" "  -0.63%  "

# Row 38
Set-CellText "D38" "1.80"
Set-CellText "E38" "  +1.30%  "

# Row 39
Set-CellText "D39" "5.52"
Set-CellText "E39" "  -3.36%  "

# Row 40
Set-CellText "D40" "0.0986"
Set-CellText "E40" "  +0.54%  "

# Row 41
Set-CellText "D41" "2.95"
Set-CellText "E41" "  +1.50%  "

# Row 42
Set-CellText "D42" "1.19"
Set-CellText "E42" "  +1.26%  "

# Row 43
Set-CellText "D43" "0.0215"
Set-CellText "E43" "  +1.42%  "

# Row 44
Set-CellText "D44" "16.77"
Set-CellText "E44" "  +4.97%  "

# Row 45
Set-CellText "D45" "1.384.99"
Set-CellText "E45" "  +1.34%  "

# Row 46
Set-CellText "D46" "90.95"
Set-CellText "E46" "  +3.53%  "

# Row 47
Set-CellText "E47" "  +1.00%  "

# Row 48
Set-CellText "D48" "7.26"
Set-CellText "E48" "  +1.49%  "

# Row 49
Set-CellText "B49" "MultiversX"
Set-CellText "C49" "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-CellText "D49" "46.70"
Set-CellText "E49" "  +6.87%  "

# Row 50
Set-CellText "B50" "MXToken"
Set-CellText "C50" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-CellText "D50" "2.83"
Set-CellText "E50" "  +0.27%  "

# Row 51
Set-CellText "E51" "  +12.15%  "
